# Issue #100 Only display enabled elements on Track screen
#
# Applies the following changes to the "Issues" sheet of the workbook:
#  - Issue #80 (client side logging, row 73): mark Status as DONE and correct
#    its Description (was erroneously "use ngx-logger", should be
#    "use node debug" like the related issue #81).
#  - Issue #81 (server side logging, row 74): mark Status as DONE.
#  - Issue #100 (row 100, "On play playlist screen should only show the
#    elements that are enabled"): mark Status as DONE.
#  - Issue #102 (row 103, "range sliders for config"): set Priority to 5.
#  - Two new issues are logged: "remove repository.service.ts" / "refactor
#    into another component" (row 104) and "use single quote style in html"
#    (row 105).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# --- Issue #80 / #81 (client + server side logging) ---------------------
$ws.Range("D73").Value = "DONE"
$ws.Range("I73").Value = "use node debug"
$ws.Range("D74").Value = "DONE"

# --- Issue #100 (Track screen only shows enabled elements) ---------------
$ws.Range("D100").Value = "DONE"

# --- Issue #102 (range sliders for config) gets a priority ---------------
$ws.Range("C103").Value = 5

# --- New issue: remove repository.service.ts / refactor into another component
$ws.Range("A104").Value = 104
$ws.Range("C104").Value = 3
$ws.Range("F104").Value = "remove repository.service.ts"
$ws.Range("I104").Value = "refactor into another component"
$ws.Rows.Item(104).RowHeight = 29

# --- New issue: use single quote style in html ---------------------------
$ws.Range("A105").Value = 105
$ws.Range("C105").Value = 3
$ws.Range("F105").Value = "use single quote style in html"
$ws.Rows.Item(105).RowHeight = 29

# --- Restore the active selection to D100, matching the saved view state -
$null = $ws.Range("D100").Select()
